# Update the "Hortaliza, Terminal La Palmera de La Serena - Melon" sheet with
# the latest weekly price rows. Three groups of new rows are introduced:
#   1. One new row at row 2 (most recent week for the "caja 12 unidades" box).
#   2. Three new rows after the existing row 47 (a new sampling date/variety
#      grouping: Primera / Segunda / Tercera).
#   3. The remaining (pre-existing) rows shift down accordingly; no explicit
#      action is required for the trailing rows since they are simply pushed
#      down by the two inserts above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by every data row in this sheet.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112027
$categoria = "Melón"
$clasif    = "Hortaliza"

# ---------------------------------------------------------------------
# 1) Insert the new first data row (row 2).
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = $mercadoId
$ws.Range("B2").Value = $mercado
$ws.Range("C2").Value = $region
$ws.Range("D2").Value = 44529
$ws.Range("E2").Value = $codreg
$ws.Range("F2").Value = $catId
$ws.Range("G2").Value = $categoria
$ws.Range("H2").Value = "Tuna"
$ws.Range("I2").Value = "Extra"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 23000
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 23500
$ws.Range("N2").Value = "$/caja 12 unidades"
$ws.Range("O2").Value = "Provincia de Copiapó"
$ws.Range("P2").Value = 1958
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = $clasif

# ---------------------------------------------------------------------
# 2) Insert three new rows after the current row 47 (before old row 48).
# ---------------------------------------------------------------------
$ws.Range("A48:A50").EntireRow.Insert()
$ws.Range("A48:R50").ClearFormats()
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 48
$ws.Range("A48").Value = $mercadoId
$ws.Range("B48").Value = $mercado
$ws.Range("C48").Value = $region
$ws.Range("D48").Value = 44904
$ws.Range("E48").Value = $codreg
$ws.Range("F48").Value = $catId
$ws.Range("G48").Value = $categoria
$ws.Range("H48").Value = "Tuna"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 1400
$ws.Range("K48").Value = 1900
$ws.Range("L48").Value = 2000
$ws.Range("M48").Value = 1950
$ws.Range("N48").Value = "$/unidad"
$ws.Range("O48").Value = "Región de O'Higgins"
$ws.Range("P48").Value = 1950
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = $clasif

# Row 49
$ws.Range("A49").Value = $mercadoId
$ws.Range("B49").Value = $mercado
$ws.Range("C49").Value = $region
$ws.Range("D49").Value = 44904
$ws.Range("E49").Value = $codreg
$ws.Range("F49").Value = $catId
$ws.Range("G49").Value = $categoria
$ws.Range("H49").Value = "Tuna"
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 1600
$ws.Range("M49").Value = 1550
$ws.Range("N49").Value = "$/unidad"
$ws.Range("O49").Value = "Región de O'Higgins"
$ws.Range("P49").Value = 1550
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = $clasif

# Row 50
$ws.Range("A50").Value = $mercadoId
$ws.Range("B50").Value = $mercado
$ws.Range("C50").Value = $region
$ws.Range("D50").Value = 44904
$ws.Range("E50").Value = $codreg
$ws.Range("F50").Value = $catId
$ws.Range("G50").Value = $categoria
$ws.Range("H50").Value = "Tuna"
$ws.Range("I50").Value = "Tercera"
$ws.Range("J50").Value = 800
$ws.Range("K50").Value = 1200
$ws.Range("L50").Value = 1300
$ws.Range("M50").Value = 1250
$ws.Range("N50").Value = "$/unidad"
$ws.Range("O50").Value = "Región de O'Higgins"
$ws.Range("P50").Value = 1250
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = $clasif
